$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, pushing existing row 120 (and below) down to 121.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new daily price observation.
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44587
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = 100112023
$ws.Cells.Item(120, 7).Value = "Brócoli"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 300
$ws.Cells.Item(120, 11).Value = 700
$ws.Cells.Item(120, 12).Value = 750
$ws.Cells.Item(120, 13).Value = 725
$ws.Cells.Item(120, 14).Value = "`$/unidad"
$ws.Cells.Item(120, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(120, 16).Value = 725
$ws.Cells.Item(120, 17).Value = 1
$ws.Cells.Item(120, 18).Value = "Hortaliza"
